$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 81.05837566666666
$ws.Range("H2").Value = 243.175127
$ws.Range("I2").Value = 0.3545816884225585
$ws.Range("J2").Value = 0.3545816884225585
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 31.22896466666667
$ws.Range("N2").Value = 93.686894
$ws.Range("O2").Value = 0.2877106972998646
$ws.Range("P2").Value = 0.2877106972998646
$ws.Range("Q2").Value = 2531.369149631726
$ws.Range("R2").Value = 22782.32234668553
$ws.Range("S2").Value = 0.1020169448258176
$ws.Range("T2").Value = 0.1020169448258176

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 81.05837566666666
$ws.Range("H3").Value = 243.175127
$ws.Range("I3").Value = 0.3545816884225585
$ws.Range("J3").Value = 0.3545816884225585
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 40.44578266666667
$ws.Range("N3").Value = 121.337348
$ws.Range("O3").Value = 0.3726247238124506
$ws.Range("P3").Value = 0.3726247238124505
$ws.Range("Q3").Value = 3278.469445527021
$ws.Range("R3").Value = 29506.22500974319
$ws.Range("S3").Value = 0.1321259037174083
$ws.Range("T3").Value = 0.1321259037174083

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 81.05837566666666
$ws.Range("H4").Value = 243.175127
$ws.Range("I4").Value = 0.3545816884225585
$ws.Range("J4").Value = 0.3545816884225585
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 25.36964133333333
$ws.Range("N4").Value = 76.108924
$ws.Range("O4").Value = 0.2337290805561598
$ws.Range("P4").Value = 0.2337290805561598
$ws.Range("Q4").Value = 2056.421917725927
$ws.Range("R4").Value = 18507.79725953334
$ws.Range("S4").Value = 0.08287605201705532
$ws.Range("T4").Value = 0.08287605201705532

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 81.05837566666666
$ws.Range("H5").Value = 243.175127
$ws.Range("I5").Value = 0.3545816884225585
$ws.Range("J5").Value = 0.3545816884225585
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 11.49855033333333
$ws.Range("N5").Value = 34.495651
$ws.Range("O5").Value = 0.1059354983315251
$ws.Range("P5").Value = 0.1059354983315251
$ws.Range("Q5").Value = 932.0538125414085
$ws.Range("R5").Value = 8388.484312872677
$ws.Range("S5").Value = 0.03756278786227731
$ws.Range("T5").Value = 0.03756278786227731

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 123.018252
$ws.Range("H6").Value = 369.054756
$ws.Range("I6").Value = 0.5381309351710768
$ws.Range("J6").Value = 0.5381309351710768
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 31.22896466666667
$ws.Range("N6").Value = 93.686894
$ws.Range("O6").Value = 0.2877106972998646
$ws.Range("P6").Value = 0.2877106972998646
$ws.Range("Q6").Value = 3841.732645063096
$ws.Range("R6").Value = 34575.59380556786
$ws.Range("S6").Value = 0.1548260265966987
$ws.Range("T6").Value = 0.1548260265966987

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 123.018252
$ws.Range("H7").Value = 369.054756
$ws.Range("I7").Value = 0.5381309351710768
$ws.Range("J7").Value = 0.5381309351710768
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 40.44578266666667
$ws.Range("N7").Value = 121.337348
$ws.Range("O7").Value = 0.3726247238124506
$ws.Range("P7").Value = 0.3726247238124505
$ws.Range("Q7").Value = 4975.569484425232
$ws.Range("R7").Value = 44780.12535982709
$ws.Range("S7").Value = 0.2005208910930583
$ws.Range("T7").Value = 0.2005208910930582

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 123.018252
$ws.Range("H8").Value = 369.054756
$ws.Range("I8").Value = 0.5381309351710768
$ws.Range("J8").Value = 0.5381309351710768
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 25.36964133333333
$ws.Range("N8").Value = 76.108924
$ws.Range("O8").Value = 0.2337290805561598
$ws.Range("P8").Value = 0.2337290805561598
$ws.Range("Q8").Value = 3120.928930693616
$ws.Range("R8").Value = 28088.36037624254
$ws.Range("S8").Value = 0.1257768486963622
$ws.Range("T8").Value = 0.1257768486963622

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 123.018252
$ws.Range("H9").Value = 369.054756
$ws.Range("I9").Value = 0.5381309351710768
$ws.Range("J9").Value = 0.5381309351710768
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 11.49855033333333
$ws.Range("N9").Value = 34.495651
$ws.Range("O9").Value = 0.1059354983315251
$ws.Range("P9").Value = 0.1059354983315251
$ws.Range("Q9").Value = 1414.531562540684
$ws.Range("R9").Value = 12730.78406286616
$ws.Range("S9").Value = 0.05700716878495767
$ws.Range("T9").Value = 0.05700716878495767

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.3624666666666667
$ws.Range("H10").Value = 1.0874
$ws.Range("I10").Value = 0.001585573873230423
$ws.Range("J10").Value = 0.001585573873230423
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 31.22896466666667
$ws.Range("N10").Value = 93.686894
$ws.Range("O10").Value = 0.2877106972998646
$ws.Range("P10").Value = 0.2877106972998646
$ws.Range("Q10").Value = 11.31945872617778
$ws.Range("R10").Value = 101.8751285356
$ws.Range("S10").Value = 0.0004561865646875722
$ws.Range("T10").Value = 0.0004561865646875722

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.3624666666666667
$ws.Range("H11").Value = 1.0874
$ws.Range("I11").Value = 0.001585573873230423
$ws.Range("J11").Value = 0.001585573873230423
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 40.44578266666667
$ws.Range("N11").Value = 121.337348
$ws.Range("O11").Value = 0.3726247238124506
$ws.Range("P11").Value = 0.3726247238124505
$ws.Range("Q11").Value = 14.66024802391111
$ws.Range("R11").Value = 131.9422322152
$ws.Range("S11").Value = 0.0005908240265967241
$ws.Range("T11").Value = 0.000590824026596724

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.3624666666666667
$ws.Range("H12").Value = 1.0874
$ws.Range("I12").Value = 0.001585573873230423
$ws.Range("J12").Value = 0.001585573873230423
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 25.36964133333333
$ws.Range("N12").Value = 76.108924
$ws.Range("O12").Value = 0.2337290805561598
$ws.Range("P12").Value = 0.2337290805561598
$ws.Range("Q12").Value = 9.195649328622222
$ws.Range("R12").Value = 82.7608439576
$ws.Range("S12").Value = 0.0003705947235440158
$ws.Range("T12").Value = 0.0003705947235440158

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.3624666666666667
$ws.Range("H13").Value = 1.0874
$ws.Range("I13").Value = 0.001585573873230423
$ws.Range("J13").Value = 0.001585573873230423
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 11.49855033333333
$ws.Range("N13").Value = 34.495651
$ws.Range("O13").Value = 0.1059354983315251
$ws.Range("P13").Value = 0.1059354983315251
$ws.Range("Q13").Value = 4.167841210822222
$ws.Range("R13").Value = 37.5105708974
$ws.Range("S13").Value = 0.0001679685584021114
$ws.Range("T13").Value = 0.0001679685584021114

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 24.16373066666667
$ws.Range("H14").Value = 72.491192
$ws.Range("I14").Value = 0.1057018025331343
$ws.Range("J14").Value = 0.1057018025331344
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 31.22896466666667
$ws.Range("N14").Value = 93.686894
$ws.Range("O14").Value = 0.2877106972998646
$ws.Range("P14").Value = 0.2877106972998646
$ws.Range("Q14").Value = 754.608291204183
$ws.Range("R14").Value = 6791.474620837647
$ws.Range("S14").Value = 0.03041153931266067
$ws.Range("T14").Value = 0.03041153931266068

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 24.16373066666667
$ws.Range("H15").Value = 72.491192
$ws.Range("I15").Value = 0.1057018025331343
$ws.Range("J15").Value = 0.1057018025331344
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 40.44578266666667
$ws.Range("N15").Value = 121.337348
$ws.Range("O15").Value = 0.3726247238124506
$ws.Range("P15").Value = 0.3726247238124505
$ws.Range("Q15").Value = 977.3209989598685
$ws.Range("R15").Value = 8795.888990638816
$ws.Range("S15").Value = 0.03938710497538737
$ws.Range("T15").Value = 0.03938710497538737

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 24.16373066666667
$ws.Range("H16").Value = 72.491192
$ws.Range("I16").Value = 0.1057018025331343
$ws.Range("J16").Value = 0.1057018025331344
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 25.36964133333333
$ws.Range("N16").Value = 76.108924
$ws.Range("O16").Value = 0.2337290805561598
$ws.Range("P16").Value = 0.2337290805561598
$ws.Range("Q16").Value = 613.0251802886008
$ws.Range("R16").Value = 5517.226622597408
$ws.Range("S16").Value = 0.02470558511919825
$ws.Range("T16").Value = 0.02470558511919825

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 24.16373066666667
$ws.Range("H17").Value = 72.491192
$ws.Range("I17").Value = 0.1057018025331343
$ws.Range("J17").Value = 0.1057018025331344
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 11.49855033333333
$ws.Range("N17").Value = 34.495651
$ws.Range("O17").Value = 0.1059354983315251
$ws.Range("P17").Value = 0.1059354983315251
$ws.Range("Q17").Value = 277.8478733117769
$ws.Range("R17").Value = 2500.630859805992
$ws.Range("S17").Value = 0.01119757312588805
$ws.Range("T17").Value = 0.01119757312588805

Write-Output "Updated cells for rows 2-17"